$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5
# from 2023-11-03 (serial 45233) to 2023-11-13 (serial 45243)
$ws.Range("C2:C5").Value = 45243
